# Update cryptocurrency price/volume data per the scraped GitHub Actions commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.730.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.777.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.773.79'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.513'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  -2.96%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.63'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.57%  '
$ws.Range('E13').Value = '  -3.54%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.407.27'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.759.61'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.711.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '456.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.36%  '
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('E25').Value = '  -6.19%  '
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('E29').Value = '  -1.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.918.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.58'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.83%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.85'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0988'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.144'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.76'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.978'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.17'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.61%  '
$ws.Range('E47').Value = '  -2.20%  '
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '384.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.31%  '
